# Applies the trade-log update described in the commit:
# "Trade #118 closed at 2026-02-16 21:44:36 - leadlag UP +0.000%"
#
# - Trade #99 (leadlag row 75 / All Trades row 100) transitions from
#   OPEN to CLOSED with exit price / pnl / exit-reason / duration filled in.
# - Trade #118 is newly logged as OPEN, appended to both the "leadlag"
#   sheet (row 93) and the "All Trades" sheet (row 100).
# - The Summary and Comparison roll-up sheets are refreshed to reflect
#   the now 91 leadlag trades (99 total trades overall).

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($Cell, [string]$Text)
    # Force a literal-text interpretation so Excel does not silently
    # reinterpret dates / percentages / numeric-looking strings.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
}

# ---------------------------------------------------------------------
# Sheet "leadlag": update trade #99 (row 75) from OPEN -> CLOSED
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

Set-TextCell $leadlag.Cells.Item(75, 8) "CLOSED"      # H75 Status
$leadlag.Cells.Item(75, 7).Value  = 67836.181931      # G75 Exit Price
$leadlag.Cells.Item(75, 9).Value  = 0.7612            # I75 P&L %
$leadlag.Cells.Item(75, 10).Value = 7.61              # J75 P&L $
Set-TextCell $leadlag.Cells.Item(75, 13) "time_exit_5min"  # M75 Exit Reason
$leadlag.Cells.Item(75, 14).Value = 5                 # N75 Duration (min)

# ---------------------------------------------------------------------
# Sheet "leadlag": append new trade #118 (row 93), status OPEN
# ---------------------------------------------------------------------
$leadlag.Cells.Item(93, 1).Value = 118                # A93 Trade #
Set-TextCell $leadlag.Cells.Item(93, 2) "2026-02-16"  # B93 Date
Set-TextCell $leadlag.Cells.Item(93, 3) "21:44:36"    # C93 Time
Set-TextCell $leadlag.Cells.Item(93, 4) "leadlag"     # D93 Strategy
Set-TextCell $leadlag.Cells.Item(93, 5) "UP"          # E93 Side
$leadlag.Cells.Item(93, 6).Value = 68456.995          # F93 Entry Price
Set-TextCell $leadlag.Cells.Item(93, 7) ""            # G93 Exit Price (empty)
Set-TextCell $leadlag.Cells.Item(93, 8) "OPEN"        # H93 Status
$leadlag.Cells.Item(93, 9).Value = 0                  # I93 P&L %
$leadlag.Cells.Item(93, 10).Value = 0                 # J93 P&L $
$leadlag.Cells.Item(93, 11).Value = 0.6098            # K93 Confidence
Set-TextCell $leadlag.Cells.Item(93, 12) "Coinbase leading with 0.061% move"  # L93 Entry Reason
Set-TextCell $leadlag.Cells.Item(93, 13) ""           # M93 Exit Reason (empty)
$leadlag.Cells.Item(93, 14).Value = 0                 # N93 Duration (min)

# ---------------------------------------------------------------------
# Sheet "All Trades": append trade #99 (row 100), now CLOSED
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(100, 1).Value = 99              # A100 Trade #
Set-TextCell $allTrades.Cells.Item(100, 2) "2026-02-16"  # B100 Date
Set-TextCell $allTrades.Cells.Item(100, 3) "21:39:31"    # C100 Time
Set-TextCell $allTrades.Cells.Item(100, 4) "leadlag"     # D100 Strategy
Set-TextCell $allTrades.Cells.Item(100, 5) "DOWN"        # E100 Side
$allTrades.Cells.Item(100, 6).Value = 68356.52499999999  # F100 Entry Price
$allTrades.Cells.Item(100, 7).Value = 67836.181931       # G100 Exit Price
Set-TextCell $allTrades.Cells.Item(100, 8) "CLOSED"      # H100 Status
$allTrades.Cells.Item(100, 9).Value  = 0.7612            # I100 P&L %
$allTrades.Cells.Item(100, 10).Value = 7.61              # J100 P&L $
$allTrades.Cells.Item(100, 11).Value = 0.75              # K100 Confidence
Set-TextCell $allTrades.Cells.Item(100, 12) "Binance leading with -0.084% move"  # L100 Entry Reason
Set-TextCell $allTrades.Cells.Item(100, 13) "time_exit_5min"  # M100 Exit Reason
$allTrades.Cells.Item(100, 14).Value = 5                 # N100 Duration (min)

# ---------------------------------------------------------------------
# Sheet "Summary": refresh OVERALL and leadlag roll-up rows
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 99
Set-TextCell $summary.Range("D2") "72.7%"
Set-TextCell $summary.Range("E2") "+31.8184%"
Set-TextCell $summary.Range("F2") "+0.3214%"

$summary.Range("C3").Value = 91
Set-TextCell $summary.Range("D3") "53.8%"
Set-TextCell $summary.Range("E3") "+18.0078%"
Set-TextCell $summary.Range("F3") "+0.1979%"

# ---------------------------------------------------------------------
# Sheet "Comparison": refresh leadlag roll-up row
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B2").Value = 91
Set-TextCell $comparison.Range("C2") "53.8%"
Set-TextCell $comparison.Range("D2") "3.35"
Set-TextCell $comparison.Range("E2") "+0.5236%"
Set-TextCell $comparison.Range("G2") "1.71"
